$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize case "unassigned" -> "Unassigned" for rows 29, 41, 60 (columns B, C, D)
$ws.Range("B29").Value = "Unassigned"
$ws.Range("C29").Value = "Unassigned"
$ws.Range("D29").Value = "Unassigned"

$ws.Range("B41").Value = "Unassigned"
$ws.Range("C41").Value = "Unassigned"
$ws.Range("D41").Value = "Unassigned"

$ws.Range("B60").Value = "Unassigned"
$ws.Range("C60").Value = "Unassigned"
$ws.Range("D60").Value = "Unassigned"

# Swap full row content (A,B,C,D) between row 42 and row 43
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "Unassigned"
$ws.Range("C42").Value = "Unassigned"
$ws.Range("D42").Value = "Unassigned"

$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"

# Swap full row content (A,B,C,D) between row 55 and row 56
$ws.Range("A55").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B55").Value = "Centropristis striata"
$ws.Range("C55").Value = "Black sea bass"
$ws.Range("D55").Value = "Teleost Fish"

$ws.Range("A56").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B56").Value = "Unassigned"
$ws.Range("C56").Value = "Unassigned"
$ws.Range("D56").Value = "Unassigned"
